# Applies the "work on data for preselection" commit:
#  - flips the sign of the month-over-month % change formulas on the
#    "IP index" sheet (C3:C388): (prev-curr)/prev*100 -> (curr-prev)/prev*100
#  - updates sheet selections (gdp growth -> B130, IP index -> D4)
#  - makes "IP index" the active/selected sheet (was "Info")

$wb = $excel.ActiveWorkbook

# --- gdp growth: just move the cursor / remembered selection ---
$wsGdp = $wb.Worksheets.Item("gdp growth")
$wsGdp.Activate() | Out-Null
$wsGdp.Range("B130").Select() | Out-Null

# --- IP index: flip the sign of every C3:C388 growth-rate formula ---
$wsIp = $wb.Worksheets.Item("IP index")
$wsIp.Activate() | Out-Null

$wsIp.Range("C3").Formula = "=(B3-B2)/B2*100"
$wsIp.Range("C4:C67").Formula = "=(B4-B3)/B3*100"
$wsIp.Range("C68:C131").Formula = "=(B68-B67)/B67*100"
$wsIp.Range("C132:C195").Formula = "=(B132-B131)/B131*100"
$wsIp.Range("C196:C259").Formula = "=(B196-B195)/B195*100"
$wsIp.Range("C260:C323").Formula = "=(B260-B259)/B259*100"
$wsIp.Range("C324:C387").Formula = "=(B324-B323)/B323*100"
$wsIp.Range("C388").Formula = "=(B388-B387)/B387*100"

# Leave the cursor where the commit's diff shows it, and make "IP index"
# the active sheet (tabSelected) last, since it's the final sheet shown.
$wsIp.Range("D4").Select() | Out-Null
